# Fruta / hortaliza, semanal
# Insert 3 new weekly-report rows for "Artic Star" right above the existing
# block that starts at row 323, pushing the existing rows (and the two
# trailing "Super Queen" rows that were mis-dated) down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 323 (shifts 323:332 down to 326:335).
$ws.Range("A323:A325").EntireRow.Insert()

# Row 323 - Comercializadora del Agro de Limarí, Nectarín, Artic Star, Especial
$ws.Range("D323").Value = 45267
$ws.Range("K323").Value = "Artic Star"
$ws.Range("L323").Value = "Especial"
$ws.Range("M323").Value = 10
$ws.Range("N323").Value = 400000
$ws.Range("O323").Value = 410000
$ws.Range("P323").Value = 405000
$ws.Range("Q323").Value = "$/bins (420 kilos)"
$ws.Range("R323").Value = "Región de O'Higgins"
$ws.Range("S323").Value = 964
$ws.Range("T323").Value = 420

# Row 324 - Artic Star, Primera
$ws.Range("D324").Value = 45267
$ws.Range("K324").Value = "Artic Star"
$ws.Range("L324").Value = "Primera"
$ws.Range("M324").Value = 16
$ws.Range("N324").Value = 370000
$ws.Range("O324").Value = 380000
$ws.Range("P324").Value = 375000
$ws.Range("Q324").Value = "$/bins (420 kilos)"
$ws.Range("R324").Value = "Región de O'Higgins"
$ws.Range("S324").Value = 893
$ws.Range("T324").Value = 420

# Row 325 - Artic Star, Segunda
$ws.Range("D325").Value = 45267
$ws.Range("K325").Value = "Artic Star"
$ws.Range("L325").Value = "Segunda"
$ws.Range("M325").Value = 16
$ws.Range("N325").Value = 330000
$ws.Range("O325").Value = 340000
$ws.Range("P325").Value = 335000
$ws.Range("Q325").Value = "$/bins (420 kilos)"
$ws.Range("R325").Value = "Región de O'Higgins"
$ws.Range("S325").Value = 798
$ws.Range("T325").Value = 420

# Fill in the remaining columns (A,B,C,E,F,G,H,I,J) that share the same
# constant values across every data row in this sheet.
$ws.Range("A323:A325").Value = 2
$ws.Range("B323:B325").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C323:C325").Value = "Coquimbo"
$ws.Range("E323:E325").Value = 4
$ws.Range("F323:F325").Value = "Fruta"
$ws.Range("G323:G325").Value = 100103
$ws.Range("H323:H325").Value = "Frutos de hueso (carozo)"
$ws.Range("I323:I325").Value = 100103006
$ws.Range("J323:J325").Value = "Nectarín"
